$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4348
$ws.Range("I40").Value = 2996.6667
$ws.Range("J40").Value = 4586.4707
$ws.Range("K40").Value = 2996.6667
$ws.Range("L40").Value = 4586.4707
$ws.Range("M40").Value = -2821.6667
$ws.Range("N40").Value = -4936.4707
$ws.Range("H58").Value = 14709938
$ws.Range("I58").Value = 22727768
$ws.Range("J58").Value = 10583.333
$ws.Range("K58").Value = 68183304
$ws.Range("L58").Value = 31749.999
$ws.Range("M58").Value = -68183154
$ws.Range("N58").Value = -32049.999
$ws.Range("H98").Value = 9159.474
$ws.Range("I98").Value = 12948.462
$ws.Range("K98").Value = 12948.462
$ws.Range("M98").Value = -11450.462
$ws.Range("H113").Value = 3774.6843
$ws.Range("I113").Value = 3897.5
$ws.Range("J113").Value = 3119.6667
$ws.Range("K113").Value = 3897.5
$ws.Range("L113").Value = 3119.6667
$ws.Range("M113").Value = -643.5
$ws.Range("N113").Value = -9627.6667
$ws.Range("H122").Value = 9159.474
$ws.Range("I122").Value = 12948.462
$ws.Range("K122").Value = 38845.386
$ws.Range("M122").Value = -36395.386
$ws.Range("H125").Value = 40951.2
$ws.Range("I125").Value = 1086
$ws.Range("K125").Value = 9774
$ws.Range("M125").Value = -7314
$ws.Range("H129").Value = 1241.6364
$ws.Range("J129").Value = 3000
$ws.Range("L129").Value = 9000
$ws.Range("N129").Value = -19000
$ws.Range("H138").Value = 1723.4791
$ws.Range("J138").Value = 1810.3784
$ws.Range("L138").Value = 5431.135200000001
$ws.Range("N138").Value = -15711.1352
$ws.Range("H141").Value = 3559.4119
$ws.Range("I141").Value = 3286.5715
$ws.Range("J141").Value = 4832.6665
$ws.Range("K141").Value = 9859.7145
$ws.Range("L141").Value = 14497.9995
$ws.Range("M141").Value = -4679.7145
$ws.Range("N141").Value = -24857.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3014
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H63").Value = 9521.299999999999
$ws.Range("I63").Value = 8498.799999999999
$ws.Range("K63").Value = 8498.799999999999
$ws.Range("M63").Value = -7812.799999999999
$ws.Range("H66").Value = 9521.299999999999
$ws.Range("I66").Value = 8498.799999999999
$ws.Range("K66").Value = 42494
$ws.Range("M66").Value = -39062
$ws.Range("H74").Value = 1998.2858
$ws.Range("I74").Value = 2040.9445
$ws.Range("J74").Value = 1742.3334
$ws.Range("K74").Value = 2040.9445
$ws.Range("L74").Value = 1742.3334
$ws.Range("M74").Value = -1166.9445
$ws.Range("N74").Value = -3490.3334
$ws.Range("H77").Value = 1998.2858
$ws.Range("I77").Value = 2040.9445
$ws.Range("J77").Value = 1742.3334
$ws.Range("K77").Value = 10204.7225
$ws.Range("L77").Value = 8711.666999999999
$ws.Range("M77").Value = -5836.7225
$ws.Range("N77").Value = -17447.667
$ws.Range("H110").Value = 877.5
$ws.Range("I110").Value = 755.5
$ws.Range("J110").Value = 999.5
$ws.Range("K110").Value = 755.5
$ws.Range("L110").Value = 999.5
$ws.Range("M110").Value = 1289.5
$ws.Range("N110").Value = -5089.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H94").Value = 3182
$ws.Range("I94").Value = 2208.0715
$ws.Range("J94").Value = 9999.5
$ws.Range("K94").Value = 2208.0715
$ws.Range("L94").Value = 9999.5
$ws.Range("M94").Value = -1757.0715
$ws.Range("N94").Value = -10901.5
$ws.Range("H105").Value = 3201.8333
$ws.Range("I105").Value = 2682.28
$ws.Range("J105").Value = 5799.6
$ws.Range("K105").Value = 2682.28
$ws.Range("L105").Value = 5799.6
$ws.Range("M105").Value = -935.2800000000002
$ws.Range("N105").Value = -9293.6
$ws.Range("H107").Value = 3388.7
$ws.Range("I107").Value = 2123.05
$ws.Range("J107").Value = 5920
$ws.Range("K107").Value = 2123.05
$ws.Range("L107").Value = 5920
$ws.Range("M107").Value = -203.0500000000002
$ws.Range("N107").Value = -9760

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6503.5
$ws.Range("I16").Value = 4505.5
$ws.Range("K16").Value = 4505.5
$ws.Range("M16").Value = -4218.5
$ws.Range("H50").Value = 45000
$ws.Range("J50").Value = 45000
$ws.Range("L50").Value = 45000
$ws.Range("N50").Value = -46250
$ws.Range("H68").Value = 67500
$ws.Range("J68").Value = 67500
$ws.Range("L68").Value = 67500
$ws.Range("N68").Value = -68998
$ws.Range("H71").Value = 67500
$ws.Range("J71").Value = 67500
$ws.Range("L71").Value = 202500
$ws.Range("N71").Value = -209988
$ws.Range("H105").Value = 2441.9285
$ws.Range("J105").Value = 2423.1667
$ws.Range("L105").Value = 2423.1667
$ws.Range("N105").Value = -5917.1667
$ws.Range("H113").Value = 6503.5
$ws.Range("I113").Value = 4505.5
$ws.Range("K113").Value = 4505.5
$ws.Range("M113").Value = -2335.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M12").Value = 146
$ws.Range("N12").ClearContents()
$ws.Range("H34").Value = 467.5
$ws.Range("J34").Value = 471
$ws.Range("L34").Value = 1413
$ws.Range("N34").Value = -1581
$ws.Range("H68").Value = 1001.5
$ws.Range("I68").Value = 668.6667
$ws.Range("K68").Value = 2006.0001
$ws.Range("M68").Value = -1195.0001
$ws.Range("H71").Value = 1001.5
$ws.Range("I71").Value = 668.6667
$ws.Range("K71").Value = 6018.0003
$ws.Range("M71").Value = -1962.0003
$ws.Range("H87").Value = 5821
$ws.Range("I87").Value = 5821
$ws.Range("K87").Value = 17463
$ws.Range("M87").Value = -16215
$ws.Range("H90").Value = 5821
$ws.Range("I90").Value = 5821
$ws.Range("K90").Value = 52389
$ws.Range("M90").Value = -46149
$ws.Range("H103").Value = 1140.5
$ws.Range("I103").Value = 248.33333
$ws.Range("J103").Value = 2032.6666
$ws.Range("K103").Value = 744.99999
$ws.Range("L103").Value = 6097.9998
$ws.Range("M103").Value = 134.00001
$ws.Range("N103").Value = -7855.9998
$ws.Range("H139").Value = 7083.3335
$ws.Range("I139").Value = 30000
$ws.Range("K139").Value = 90000
$ws.Range("M139").Value = -84860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 68333
$ws.Range("J135").Value = 68333
$ws.Range("L135").Value = 68333
$ws.Range("N135").Value = -78473

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6341
$ws.Range("I7").Value = 5714
$ws.Range("J7").Value = 7218.8
$ws.Range("K7").Value = 5714
$ws.Range("L7").Value = 7218.8
$ws.Range("M7").Value = -5602
$ws.Range("N7").Value = -7442.8
$ws.Range("H40").Value = 5873.0527
$ws.Range("I40").Value = 4162.5454
$ws.Range("K40").Value = 4162.5454
$ws.Range("M40").Value = -4026.5454
$ws.Range("H64").Value = 17074.5
$ws.Range("J64").Value = 17074.5
$ws.Range("L64").Value = 17074.5
$ws.Range("N64").Value = -17524.5
$ws.Range("H67").Value = 17074.5
$ws.Range("J67").Value = 17074.5
$ws.Range("L67").Value = 17074.5
$ws.Range("N67").Value = -18634.5
$ws.Range("H70").Value = 40000
$ws.Range("I70").Value = 45000
$ws.Range("J70").Value = 35000
$ws.Range("K70").Value = 45000
$ws.Range("L70").Value = 35000
$ws.Range("M70").Value = -44730
$ws.Range("N70").Value = -35540
$ws.Range("H73").Value = 40000
$ws.Range("I73").Value = 45000
$ws.Range("J73").Value = 35000
$ws.Range("K73").Value = 45000
$ws.Range("L73").Value = 35000
$ws.Range("M73").Value = -44064
$ws.Range("N73").Value = -36872
$ws.Range("H126").Value = 6341
$ws.Range("I126").Value = 5714
$ws.Range("J126").Value = 7218.8
$ws.Range("K126").Value = 17142
$ws.Range("L126").Value = 21656.4
$ws.Range("M126").Value = -14672
$ws.Range("N126").Value = -26596.4
